$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (row 47),
# shifting the existing data (old rows 47-145) down to rows 49-147.
$ws.Rows("47:48").Insert()

# Populate the two newly inserted rows with the new price-report entries.
$ws.Cells.Item(47, 1).Value = 7
$ws.Cells.Item(47, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(47, 3).Value = 'Ñuble'
$ws.Cells.Item(47, 4).Value = 44497
$ws.Cells.Item(47, 5).Value = 16
$ws.Cells.Item(47, 6).Value = 'Fruta'
$ws.Cells.Item(47, 7).Value = 100101
$ws.Cells.Item(47, 8).Value = 'Berries'
$ws.Cells.Item(47, 9).Value = 100112025
$ws.Cells.Item(47, 10).Value = 'Frutilla'
$ws.Cells.Item(47, 11).Value = 'Sin especificar'
$ws.Cells.Item(47, 12).Value = 'Primera'
$ws.Cells.Item(47, 13).Value = 160
$ws.Cells.Item(47, 14).Value = 7000
$ws.Cells.Item(47, 15).Value = 7500
$ws.Cells.Item(47, 16).Value = 7250
$ws.Cells.Item(47, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(47, 18).Value = 'Provincia de Diguillín'
$ws.Cells.Item(47, 19).Value = 1036
$ws.Cells.Item(47, 20).Value = 7
$ws.Cells.Item(48, 1).Value = 7
$ws.Cells.Item(48, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(48, 3).Value = 'Ñuble'
$ws.Cells.Item(48, 4).Value = 44497
$ws.Cells.Item(48, 5).Value = 16
$ws.Cells.Item(48, 6).Value = 'Fruta'
$ws.Cells.Item(48, 7).Value = 100101
$ws.Cells.Item(48, 8).Value = 'Berries'
$ws.Cells.Item(48, 9).Value = 100112025
$ws.Cells.Item(48, 10).Value = 'Frutilla'
$ws.Cells.Item(48, 11).Value = 'Sin especificar'
$ws.Cells.Item(48, 12).Value = 'Segunda'
$ws.Cells.Item(48, 13).Value = 160
$ws.Cells.Item(48, 14).Value = 6000
$ws.Cells.Item(48, 15).Value = 6500
$ws.Cells.Item(48, 16).Value = 6250
$ws.Cells.Item(48, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(48, 18).Value = 'Provincia de Diguillín'
$ws.Cells.Item(48, 19).Value = 893
$ws.Cells.Item(48, 20).Value = 7
